$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Sentiment column (K) was re-tagged from "Negative" to "Neutral"
# for these rows (rows 4, 7, 15, 16 were left untouched).
$rows = @(2, 3, 5, 6, 8, 9, 10, 11, 12, 13, 14)
foreach ($r in $rows) {
    $ws.Range("K$r").Value = "Neutral"
}

# The active selection moved from AA8 (with the view scrolled to column R)
# back to C9 with the view scrolled to the top-left (A1).
$ws.Range("C9").Select()
